$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix existing Temp[c] readings (H column)
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 2

# Add userid batch / station id rows: pressure G9 cleared, new Temp[c]/Delta[c]
# values added for rows 9-11 (H and I columns)
$ws.Range("G9").ClearContents()

$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 1

$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 1

$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 1

# Update the active selection to match the author's final cursor position
$ws.Range("I12").Select()
